$wb = $excel.ActiveWorkbook

# --- snapshot sheet: remove returned players, refresh scraped_at timestamps ---
$snap = $wb.Worksheets.Item("snapshot")

# Remove row 19 (НХК / Дергачёв Александр - returned) and row 3 (АВТ / Кизимов Семён - returned)
# Delete the higher-indexed row first so the lower row index stays valid.
$snap.Rows.Item(19).Delete()
$snap.Rows.Item(3).Delete()

# Refresh the scraped_at (column K) timestamps for the new scrape pass
$snap.Range("K2").Value = "2025-11-21T03:01:40.482673+00:00"
$snap.Range("K3").Value = "2025-11-21T03:01:40.482705+00:00"
$snap.Range("K4").Value = "2025-11-21T03:01:42.726478+00:00"
$snap.Range("K5").Value = "2025-11-21T03:01:42.726507+00:00"
$snap.Range("K6").Value = "2025-11-21T03:01:45.074404+00:00"
$snap.Range("K7").Value = "2025-11-21T03:01:47.699807+00:00"
$snap.Range("K8").Value = "2025-11-21T03:01:50.562316+00:00"
$snap.Range("K9").Value = "2025-11-21T03:01:50.562346+00:00"
$snap.Range("K10").Value = "2025-11-21T03:01:50.562374+00:00"
$snap.Range("K11").Value = "2025-11-21T03:01:52.804251+00:00"
$snap.Range("K12").Value = "2025-11-21T03:01:55.026165+00:00"
$snap.Range("K13").Value = "2025-11-21T03:01:57.673962+00:00"
$snap.Range("K14").Value = "2025-11-21T03:01:59.788791+00:00"
$snap.Range("K15").Value = "2025-11-21T03:01:59.788823+00:00"
$snap.Range("K16").Value = "2025-11-21T03:01:59.788842+00:00"
$snap.Range("K17").Value = "2025-11-21T03:02:02.126742+00:00"
$snap.Range("K18").Value = "2025-11-21T03:02:07.440766+00:00"
$snap.Range("K19").Value = "2025-11-21T03:02:07.440794+00:00"
$snap.Range("K20").Value = "2025-11-21T03:02:10.143592+00:00"
$snap.Range("K21").Value = "2025-11-21T03:02:10.143621+00:00"
$snap.Range("K22").Value = "2025-11-21T03:02:10.143638+00:00"
$snap.Range("K23").Value = "2025-11-21T03:02:12.836980+00:00"
$snap.Range("K24").Value = "2025-11-21T03:02:12.837009+00:00"
$snap.Range("K25").Value = "2025-11-21T03:02:15.040992+00:00"
$snap.Range("K26").Value = "2025-11-21T03:02:15.041018+00:00"
$snap.Range("K27").Value = "2025-11-21T03:02:15.041037+00:00"
$snap.Range("K28").Value = "2025-11-21T03:02:17.736298+00:00"
$snap.Range("K29").Value = "2025-11-21T03:02:17.736327+00:00"
$snap.Range("K30").Value = "2025-11-21T03:02:20.020564+00:00"
$snap.Range("K31").Value = "2025-11-21T03:02:20.020628+00:00"
$snap.Range("K32").Value = "2025-11-21T03:02:20.020679+00:00"
$snap.Range("K33").Value = "2025-11-21T03:02:20.020712+00:00"
$snap.Range("K34").Value = "2025-11-21T03:02:20.020732+00:00"
$snap.Range("K35").Value = "2025-11-21T03:02:22.623304+00:00"
$snap.Range("K36").Value = "2025-11-21T03:02:22.623331+00:00"
$snap.Range("K37").Value = "2025-11-21T03:02:27.657792+00:00"
$snap.Range("K38").Value = "2025-11-21T03:02:27.657820+00:00"
$snap.Range("K39").Value = "2025-11-21T03:02:27.657837+00:00"
$snap.Range("K40").Value = "2025-11-21T03:02:30.356017+00:00"

# --- returned sheet: append the two players who returned from injury ---
$ret = $wb.Worksheets.Item("returned")
$ret.Range("A2").Value = "АВТ"
$ret.Range("B2").Value = "Автомобилист"
$ret.Range("C2").Value = "Кизимов Семён"
$ret.Range("D2").Value = "1369_АВТ_кизимовсемен"
$ret.Range("E2").Value = "RETURN"
$ret.Range("F2").Value = "2025-11-21T11:02:30.862605+08:00"
$ret.Range("G2").NumberFormat = "@"
$ret.Range("G2").Value = "2025-11-21"

$ret.Range("A3").Value = "НХК"
$ret.Range("B3").Value = "Нефтехимик"
$ret.Range("C3").Value = "Дергачёв Александр"
$ret.Range("D3").Value = "1369_НХК_дергачевалександр"
$ret.Range("E3").Value = "RETURN"
$ret.Range("F3").Value = "2025-11-21T11:02:30.862605+08:00"
$ret.Range("G3").NumberFormat = "@"
$ret.Range("G3").Value = "2025-11-21"
